# Apply updated TPM-derived values to Dhh-Hhip LR-pairs sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 4.573795666666666
$ws.Cells.Item(2, 8).Value = 13.721387
$ws.Cells.Item(2, 9).Value = 0.6529099782872276
$ws.Cells.Item(2, 10).Value = 0.6529099782872277
$ws.Cells.Item(2, 15).Value = 0.9272790495339911
$ws.Cells.Item(2, 16).Value = 0.9272790495339911
$ws.Cells.Item(2, 17).Value = 0.9716922171934445
$ws.Cells.Item(2, 18).Value = 8.745229954741
$ws.Cells.Item(2, 19).Value = 0.6054297440974392
$ws.Cells.Item(2, 20).Value = 0.6054297440974393

# Row 3
$ws.Cells.Item(3, 7).Value = 4.573795666666666
$ws.Cells.Item(3, 8).Value = 13.721387
$ws.Cells.Item(3, 9).Value = 0.6529099782872276
$ws.Cells.Item(3, 10).Value = 0.6529099782872277
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.016661
$ws.Cells.Item(3, 14).Value = 0.049983
$ws.Cells.Item(3, 15).Value = 0.07272095046600884
$ws.Cells.Item(3, 16).Value = 0.07272095046600885
$ws.Cells.Item(3, 17).Value = 0.07620400960233333
$ws.Cells.Item(3, 18).Value = 0.685836086421
$ws.Cells.Item(3, 19).Value = 0.04748023418978838
$ws.Cells.Item(3, 20).Value = 0.0474802341897884

# Row 4
$ws.Cells.Item(4, 9).Value = 0.1985019229157801
$ws.Cells.Item(4, 10).Value = 0.1985019229157801
$ws.Cells.Item(4, 15).Value = 0.9272790495339911
$ws.Cells.Item(4, 16).Value = 0.9272790495339911
$ws.Cells.Item(4, 19).Value = 0.1840666744120142
$ws.Cells.Item(4, 20).Value = 0.1840666744120142

# Row 5
$ws.Cells.Item(5, 9).Value = 0.1985019229157801
$ws.Cells.Item(5, 10).Value = 0.1985019229157801
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.016661
$ws.Cells.Item(5, 14).Value = 0.049983
$ws.Cells.Item(5, 15).Value = 0.07272095046600884
$ws.Cells.Item(5, 16).Value = 0.07272095046600885
$ws.Cells.Item(5, 17).Value = 0.023168036855
$ws.Cells.Item(5, 18).Value = 0.208512331695
$ws.Cells.Item(5, 19).Value = 0.01443524850376595
$ws.Cells.Item(5, 20).Value = 0.01443524850376596

# Row 6
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.3547236666666667
$ws.Cells.Item(6, 8).Value = 1.064171
$ws.Cells.Item(6, 9).Value = 0.05063685358512936
$ws.Cells.Item(6, 10).Value = 0.05063685358512936
$ws.Cells.Item(6, 15).Value = 0.9272790495339911
$ws.Cells.Item(6, 16).Value = 0.9272790495339911
$ws.Cells.Item(6, 17).Value = 0.07536021529477778
$ws.Cells.Item(6, 18).Value = 0.678241937653
$ws.Cells.Item(6, 19).Value = 0.04695449346381062
$ws.Cells.Item(6, 20).Value = 0.04695449346381062

# Row 7
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.3547236666666667
$ws.Cells.Item(7, 8).Value = 1.064171
$ws.Cells.Item(7, 9).Value = 0.05063685358512936
$ws.Cells.Item(7, 10).Value = 0.05063685358512936
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.016661
$ws.Cells.Item(7, 14).Value = 0.049983
$ws.Cells.Item(7, 15).Value = 0.07272095046600884
$ws.Cells.Item(7, 16).Value = 0.07272095046600885
$ws.Cells.Item(7, 17).Value = 0.005910051010333333
$ws.Cells.Item(7, 18).Value = 0.053190459093
$ws.Cells.Item(7, 19).Value = 0.003682360121318734
$ws.Cells.Item(7, 20).Value = 0.003682360121318735

# Row 8
$ws.Cells.Item(8, 7).Value = 0.479723
$ws.Cells.Item(8, 8).Value = 1.439169
$ws.Cells.Item(8, 9).Value = 0.06848052609708123
$ws.Cells.Item(8, 10).Value = 0.06848052609708123
$ws.Cells.Item(8, 15).Value = 0.9272790495339911
$ws.Cells.Item(8, 16).Value = 0.9272790495339911
$ws.Cells.Item(8, 17).Value = 0.1019160319963333
$ws.Cells.Item(8, 18).Value = 0.917244287967
$ws.Cells.Item(8, 19).Value = 0.06350055715088916
$ws.Cells.Item(8, 20).Value = 0.06350055715088916

# Row 9
$ws.Cells.Item(9, 7).Value = 0.479723
$ws.Cells.Item(9, 8).Value = 1.439169
$ws.Cells.Item(9, 9).Value = 0.06848052609708123
$ws.Cells.Item(9, 10).Value = 0.06848052609708123
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.016661
$ws.Cells.Item(9, 14).Value = 0.049983
$ws.Cells.Item(9, 15).Value = 0.07272095046600884
$ws.Cells.Item(9, 16).Value = 0.07272095046600885
$ws.Cells.Item(9, 17).Value = 0.007992664903
$ws.Cells.Item(9, 18).Value = 0.071933984127
$ws.Cells.Item(9, 19).Value = 0.00497996894619207
$ws.Cells.Item(9, 20).Value = 0.00497996894619207

# Row 10
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.2064496666666667
$ws.Cells.Item(10, 8).Value = 0.619349
$ws.Cells.Item(10, 9).Value = 0.02947071911478163
$ws.Cells.Item(10, 10).Value = 0.02947071911478163
$ws.Cells.Item(10, 15).Value = 0.9272790495339911
$ws.Cells.Item(10, 16).Value = 0.9272790495339911
$ws.Cells.Item(10, 17).Value = 0.04385974996744445
$ws.Cells.Item(10, 18).Value = 0.394737749707
$ws.Cells.Item(10, 19).Value = 0.02732758040983794
$ws.Cells.Item(10, 20).Value = 0.02732758040983794

# Row 11
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.2064496666666667
$ws.Cells.Item(11, 8).Value = 0.619349
$ws.Cells.Item(11, 9).Value = 0.02947071911478163
$ws.Cells.Item(11, 10).Value = 0.02947071911478163
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.016661
$ws.Cells.Item(11, 14).Value = 0.049983
$ws.Cells.Item(11, 15).Value = 0.07272095046600884
$ws.Cells.Item(11, 16).Value = 0.07272095046600885
$ws.Cells.Item(11, 17).Value = 0.003439657896333333
$ws.Cells.Item(11, 18).Value = 0.030956921067
$ws.Cells.Item(11, 19).Value = 0.002143138704943695
$ws.Cells.Item(11, 20).Value = 0.002143138704943696

